# Update computed results pulled from server for each scenario-year sheet.
$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "B2" = 681.7676493320005
        "E2" = 20287.66317618639
        "I2" = 11322.696948346
        "L2" = 40717.1573684766
        "M2" = 7413.4115776105
        "N2" = 5030.117123085909
        "O2" = 4892.012299288745
    }
    "2030" = @{
        "B2" = 3998.79212449002
        "E2" = 39274.23912147338
        "I2" = 30952.52893051627
        "L2" = 56191.96522533072
        "M2" = 15388.9961629327
        "N2" = 7449.833478787734
        "O2" = 8459.289075669971
    }
    "2035" = @{
        "A2" = 2861.961401238371
        "B2" = 6313.12160973443
        "E2" = 50465.92318872747
        "I2" = 45991.05621758628
        "L2" = 56191.96522533072
        "M2" = 18869.3373667282
        "N2" = 11987.77381092391
        "O2" = 11157.28118474881
    }
    "2040" = @{
        "A2" = 2861.961401238371
        "B2" = 6313.12160973443
        "E2" = 50465.92318872747
        "I2" = 45991.05621758628
        "L2" = 56191.96522533072
        "M2" = 18869.3373667282
        "N2" = 12092.18250713501
        "O2" = 11157.28118474881
    }
    "2045" = @{
        "A2" = 6302.873118834019
        "B2" = 6313.12160973443
        "E2" = 50465.92318872747
        "I2" = 45991.05621758628
        "L2" = 56191.96522533072
        "M2" = 18869.3373667282
        "N2" = 12627.69369170892
        "O2" = 13492.74717109747
    }
    "2050" = @{
        "A2" = 6302.873118834019
        "B2" = 6313.12160973443
        "E2" = 50465.92318872747
        "I2" = 45991.05621758628
        "L2" = 56191.96522533072
        "M2" = 18869.3373667282
        "N2" = 12627.69369170892
        "O2" = 13492.74717109747
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
